# Continuation of time-var experiment:
# fill in the missing F13 data point for the third ("p(x) = x_1") block
# and leave the selection where the user moved it next (I15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("F13").Value = 0.71826461805817998

$ws.Range("I15").Select()
